$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.487833380699158
$ws.Range("B1").Value = 3.064851522445679
$ws.Range("C1").Value = 2.440634250640869
$ws.Range("D1").Value = 2.301240682601929
$ws.Range("E1").Value = 1.971065998077393
